# DFS Growth.xlsx edit:
# - Add a "Theoretical" H column with the actual base (2.168) used by the
#   G-column growth-rate formulas, replacing the hardcoded 2.162 literal
#   with a reference to $H$2 so the base is easy to tweak.
# - Re-label the G-column header to "2.xxx^(n+1)" (generic, driven by H1's
#   "2.xxx" label) instead of the old literal "2.162^(n+1)".
# - Refill the A / F / G growth columns so Excel re-derives their shared
#   formula groupings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-label the existing growth-rate header ---------------------------
$ws.Range("G2").Value = "2.xxx^(n+1)"

# --- New "Theoretical" base value & label -------------------------------
$ws.Range("H1").Value = "2.xxx"
$ws.Range("H2").Value = 2.168

# --- Column A: row index (n), fill down in its two existing blocks ------
$ws.Range("A6:A55").Formula = "=A5+1"
$ws.Range("A56:A69").Formula = "=A55+1"
$ws.Range("A70:A83").Formula = "=A69+1"

# --- Column F: 2^n, fill down in three blocks ----------------------------
$ws.Range("F3:F34").Formula = "=2^A3"
$ws.Range("F35:F56").Formula = "=2^A35"
$ws.Range("F57:F83").Formula = "=2^A57"

# --- Column G: theoretical growth now driven off $H$2 -------------------
$ws.Range("G3").Formula = "=`$H`$2^(A3+1)"
$ws.Range("G4:G67").Formula = "=`$H`$2^(A4+1)"
$ws.Range("G68:G83").Formula = "=`$H`$2^(A68+1)"

# --- View tweaks matching the saved window state -------------------------
$ws.Range("G3:G83").Select()
$ws.Application.ActiveWindow.ScrollRow = 19
